# "Generate Report for Handback" - mark the eb4682be... row as handed back
# for both the zh-cn and de-de localization targets, filling in the
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns and refreshing the row Status (which also drives the Overview
# summary columns), plus widening a few columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

$srcFile   = "eb4682be-b860-4793-a44e-97f93e398258.md"
$srcUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/da5c08bebac047e2e419f4155944edd40fa1a2c1/e2e/eb4682be-b860-4793-a44e-97f93e398258.md"
$ffFile    = "ffffc624c887-4ef2-43b8-a337-6d417a577678.md"
$ffUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/da5c08bebac047e2e419f4155944edd40fa1a2c1/e2e/ffffc624c887-4ef2-43b8-a337-6d417a577678.md"

$zhHandbackFile = "eb4682be-b860-4793-a44e-97f93e398258.d865e2ceb4f0b4bfd878fdb82cd3a8a513be8135.zh-cn.xlf"
$deHandbackFile = "eb4682be-b860-4793-a44e-97f93e398258.d865e2ceb4f0b4bfd878fdb82cd3a8a513be8135.de-de.xlf"

$zhHandbackDate = "2016-10-26 07:57:24"
$deHandbackDate = "2016-10-26 07:57:42"

# Excel's ColumnWidth property is offset from the stored OOXML column width
# by the default gridline padding (~0.8333 chars); add that back so the
# saved width matches the intended character width.
$colPad = 0.8333333333333334

# ---------------------------------------------------------------------
# Overview sheet: refresh the per-language status text + widen the
# zh-cn / de-de status columns (E, F).
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusHandedBack
$wsOverview.Range("F2").Value = $statusHandedBack
$wsOverview.Range("E3").Value = $statusHandedBack
$wsOverview.Range("F3").Value = $statusHandedBack
$wsOverview.Range("E1").ColumnWidth = 29.9777050018311 - $colPad
$wsOverview.Range("F1").ColumnWidth = 29.9777050018311 - $colPad

# ---------------------------------------------------------------------
# Helper-style per-sheet update: zh-cn and de-de share the same column
# layout (Status=C, Latest Target File=I, Latest Handback File=J,
# Latest Handback DateTime=K).
# ---------------------------------------------------------------------
function Update-LangSheet($ws, $handbackFile, $handbackDate) {
    # Status column (also backs the Overview rollup above)
    $ws.Range("C2").Value = $statusHandedBack
    $ws.Range("C3").Value = $statusHandedBack

    # Latest Target File -> hyperlink to the (same) source markdown file
    $ws.Range("I2").Value = $srcFile
    $ws.Range("I2").Font.Underline = $true
    $ws.Range("I2").Font.Color = 15570276
    $ws.Range("I3").Value = $srcFile
    $ws.Range("I3").Font.Underline = $true
    $ws.Range("I3").Font.Color = 15570276

    # Latest Handback File
    $ws.Range("J2").Value = $handbackFile
    $ws.Range("J3").Value = $handbackFile

    # Latest Handback DateTime
    $ws.Range("K2").Value = $handbackDate
    $ws.Range("K3").Value = $handbackDate

    # Re-create the hyperlinks collection so ordering matches: A2, I2, A3, I3
    $ws.Range("A1:P3").Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $srcUrl, "", "", $srcFile)
    $ws.Hyperlinks.Add($ws.Range("I2"), $srcUrl, "", "", $srcFile)
    $ws.Hyperlinks.Add($ws.Range("A3"), $ffUrl, "", "", $ffFile)
    $ws.Hyperlinks.Add($ws.Range("I3"), $srcUrl, "", "", $srcFile)

    # Column widths: Status (C) + Latest Target File / Latest Handback File (I, J)
    $ws.Range("C1").ColumnWidth = 29.9777050018311 - $colPad
    $ws.Range("I1").ColumnWidth = 40 - $colPad
    $ws.Range("J1").ColumnWidth = 40 - $colPad
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-LangSheet $wsZhCn $zhHandbackFile $zhHandbackDate

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-LangSheet $wsDeDe $deHandbackFile $deHandbackDate
